$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "tsolera656"
$ws.Range("B11").Value = "Raimundo"
$ws.Range("C11").Value = "Guijarro"
$ws.Range("D11").Value = "pmanjon@gmail.com"
$ws.Range("E11").Value = "n^s4O8bE*z0m"
$ws.Range("F11").Value = "n^s4O8bE*z0m"
$ws.Range("G11").Value = "Válido"
